# Updated cryptos list on Sun Oct 15 06:27:33 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.962.18"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.556.88"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.70"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.780.90"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.568.28"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.72"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "26.955.57"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.70"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "0.0₃0704"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.57"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.55"
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.58"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.01"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.105"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +3.45%  "
$ws.Range("D34").Value = "1.421.63"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +9.96%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  +1.89%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.531"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.807"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.31"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.48"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "1.693.35"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.59"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("E51").Value = "  +0.70%  "
